$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the grammar-rule descriptions (column D) for the statement-block,
# while, and if/if-else rows: add type checking + proper scoping (push/pop)
# semantics, matching the "add more analysis functions" commit.

# <StmtBlock> ::= <Stmts>  -- now pops the locally-declared vars on exit
$ws.Range("D21").Value = "<SB>.IR = <S>.IR`n<SB>.returnType = <S>.returnType`n<SB>.innerVarAmount = <S>.vA`npop <SB>.innerVarAmount vars from var symbol table"

# while ( <Exprsn> ) <StmtBlock>  -- now type-checks <E> and uses L1/L2/L3 labels
$ws.Range("D32").Value = "if <E>.valType is not 'int': ERROR `n <WS>.returnType = <SB>.returnType`n<WS>.IR = 'L1:' + <E>.IR `n+ 'if(<E>.val == 1) goto L2 else L3:'`n+ 'L2' + <SB>.IR + 'goto L1' `n+ 'L3:'"

# if ( <Exprsn> ) <StmtBlock> else <StmtBlock>  -- now type-checks <E> and uses L1/L2/L3 labels
$ws.Range("D33").Value = "if <E>.valType is not 'int': ERROR `n<SB1>.rT equals to <SB2>.rT: <IS>.rT = <SB1>.rT`nelse: ERROR`n<IS>.IR = <E>.IR + 'if (<E>.val != 1) goto L1 else L2' `n+ 'L1:' + <SB1>.IR + 'goto L3' `n+ 'L2:' + <SB2>.IR `n+ 'L3:'"

# if ( <Exprsn> ) <StmtBlock>  -- now type-checks <E> and uses L1/L2 labels
$ws.Range("D34").Value = "if <E>.valType is not 'int': ERROR`n <IS>.returnType = <SB>.returnType`n<IS>.IR = <E>.IR + 'if (<E>.val != 1) goto L1 else L2:' `n+ 'L1:' + <SB>.IR + `n'L2:'"

# The if/if-else row grew an extra wrapped line of text, so it needs a little
# more room to display without the text spilling.
$ws.Rows.Item(21).RowHeight = 60.75

# Scroll the sheet so the newly-edited rows are in view and select D32.
$excel.ActiveWindow.ScrollRow = 32
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("D32").Select() | Out-Null
